$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header cells: "_old" suffix becomes "_FV2310", "_new" suffix becomes "_FV2404"
$baseHeaders = @("Segmentname", "Segmentgruppe", "Segment", "Datenelement", "Segment ID", "Code", "Qualifier", "Beschreibung", "Bedingungsausdruck", "Bedingung")

for ($i = 0; $i -lt $baseHeaders.Count; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value = ($baseHeaders[$i] + "_FV2310")
}
for ($i = 0; $i -lt $baseHeaders.Count; $i++) {
    $col = $i + 12
    $ws.Cells.Item(1, $col).Value = ($baseHeaders[$i] + "_FV2404")
}

# Turn the data range into an Excel Table (ListObject)
$rng = $ws.Range("A1:U76")
$lo = $ws.ListObjects.Add(1, $rng, $null, 1)
$lo.Name = "Table1"

# Freeze the header row
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
